$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the scraped cryptos.xlsx price/volume/hour/coin-name updates.
# Columns D (Price), E (Volume 1h %) and G (Hora) are stored as literal
# text in the source sheet, so force text format before assigning to
# avoid Excel auto-converting "308.24" / "-0.60%" into numeric values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.60%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "21"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.78%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "21"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.048"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.96%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "21"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07628"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.67%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "21"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.614"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.65%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "21"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9095"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.25%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "21"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.436"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.39%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "21"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1020"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.03%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "21"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1767"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.52%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "21"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09182"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.83%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "21"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04314"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.05%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "21"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1054"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.39%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "21"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001252"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.87%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "21"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005803"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-3.73%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "21"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.356"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.71%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "21"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.265"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.91%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "21"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.91%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "21"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.832"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-5.32%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "21"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1349"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.67%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "21"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2724"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.89%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "21"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04162"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.39%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "21"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001226"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.69%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "21"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004008"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.16%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "21"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001301"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.95%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "21"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003008"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.34%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "21"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "21"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "21"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "21"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "21"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "21"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "21"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "21"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "21"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "21"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "21"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "21"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02404"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.23%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "21"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05180"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.49%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "21"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007784"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.73%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "21"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1307"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.44%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "21"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007088"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-6.58%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "21"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001919"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.77%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "21"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007468"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-7.49%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "21"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3355"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.67%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "21"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006359"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-6.19%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "21"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.55%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "21"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.004402"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "6.67%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "21"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.005695"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "66.46%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "21"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.55%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "21"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.55%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "21"
